$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("22.3.2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

$rng = $d.Content
$rng.Start = $rng.Start
$rng.Find.Execute("22.3.2023")
$rng.Collapse(0)
$rng.InsertAfter(".")
